$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 66808
$ws.Range("E2").Value = 10316
$ws.Range("F2").Value = 10316
$ws.Range("G2").Value = 9356
$ws.Range("H2").Value = 6993
$ws.Range("I2").Value = 7004
$ws.Range("J2").Value = -11
$ws.Range("K2").Value = 83201
$ws.Range("L2").Value = 38036
$ws.Range("M2").Value = 45165
$ws.Range("N2").Value = 45110
$ws.Range("O2").Value = 55
$ws.Range("P2").Value = 619
$ws.Range("Q2").Value = 11470
$ws.Range("R2").Value = -9048
$ws.Range("S2").Value = -1959
$ws.Range("T2").Value = 9366
$ws.Range("U2").Value = 2104
$ws.Range("V2").Value = 24236
$ws.Range("W2").Value = 15.44
$ws.Range("X2").Value = 10.47
$ws.Range("Y2").Value = 16.62
$ws.Range("Z2").Value = 8.69
$ws.Range("AA2").Value = 84.22
$ws.Range("AB2").Value = 7323.3
$ws.Range("AC2").Value = 5654
$ws.Range("AD2").Value = 9.32
$ws.Range("AE2").Value = 36422
$ws.Range("AF2").Value = 1.45
$ws.Range("AH2").Value = 0.76
$ws.Range("AI2").Value = 7.07
$ws.Range("AJ2").Value = 123875069

$ws.Range("D3").Value = 64282
$ws.Range("E3").Value = 8850
$ws.Range("F3").Value = 8850
$ws.Range("G3").Value = 8394
$ws.Range("H3").Value = 6565
$ws.Range("I3").Value = 6554
$ws.Range("J3").Value = 11
$ws.Range("K3").Value = 94961
$ws.Range("L3").Value = 43554
$ws.Range("M3").Value = 51408
$ws.Range("N3").Value = 51344
$ws.Range("O3").Value = 64
$ws.Range("P3").Value = 619
$ws.Range("Q3").Value = 11137
$ws.Range("R3").Value = -21037
$ws.Range("S3").Value = 8887
$ws.Range("T3").Value = 8630
$ws.Range("U3").Value = 2507
$ws.Range("V3").Value = 28957
$ws.Range("W3").Value = 13.77
$ws.Range("X3").Value = 10.21
$ws.Range("Y3").Value = 13.59
$ws.Range("Z3").Value = 7.37
$ws.Range("AA3").Value = 84.72
$ws.Range("AB3").Value = 8278.530000000001
$ws.Range("AC3").Value = 5291
$ws.Range("AD3").Value = 8.890000000000001
$ws.Range("AE3").Value = 41455
$ws.Range("AF3").Value = 1.13
$ws.Range("AG3").Value = 400
$ws.Range("AH3").Value = 0.85
$ws.Range("AI3").Value = 7.56
$ws.Range("AJ3").Value = 123875069

$ws.Range("D4").Value = 66218
$ws.Range("E4").Value = 11032
$ws.Range("F4").Value = 11032
$ws.Range("G4").Value = 10924
$ws.Range("H4").Value = 8791
$ws.Range("I4").Value = 8729
$ws.Range("J4").Value = 62
$ws.Range("K4").Value = 96220
$ws.Range("L4").Value = 36605
$ws.Range("M4").Value = 59615
$ws.Range("N4").Value = 59494
$ws.Range("O4").Value = 121
$ws.Range("P4").Value = 619
$ws.Range("Q4").Value = 12178
$ws.Range("R4").Value = -7400
$ws.Range("S4").Value = -6497
$ws.Range("T4").Value = 7227
$ws.Range("U4").Value = 4951
$ws.Range("V4").Value = 23273
$ws.Range("W4").Value = 16.66
$ws.Range("X4").Value = 13.28
$ws.Range("Y4").Value = 15.75
$ws.Range("Z4").Value = 9.199999999999999
$ws.Range("AA4").Value = 61.4
$ws.Range("AB4").Value = 9596.83
$ws.Range("AC4").Value = 7046
$ws.Range("AD4").Value = 8.23
$ws.Range("AE4").Value = 48036
$ws.Range("AF4").Value = 1.21
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 0.6899999999999999
$ws.Range("AI4").Value = 5.68
$ws.Range("AJ4").Value = 123875069

$ws.Range("D5").Value = 68129
$ws.Range("E5").Value = 7934
$ws.Range("F5").Value = 7934
$ws.Range("G5").Value = 7430
$ws.Range("H5").Value = 6065
$ws.Range("I5").Value = 5991
$ws.Range("J5").Value = 74
$ws.Range("K5").Value = 95188
$ws.Range("L5").Value = 31453
$ws.Range("M5").Value = 63735
$ws.Range("N5").Value = 63587
$ws.Range("O5").Value = 148
$ws.Range("P5").Value = 619
$ws.Range("Q5").Value = 8517
$ws.Range("R5").Value = -3443
$ws.Range("S5").Value = -2817
$ws.Range("T5").Value = 4723
$ws.Range("U5").Value = 3794
$ws.Range("V5").Value = 19735
$ws.Range("W5").Value = 11.65
$ws.Range("X5").Value = 8.9
$ws.Range("Y5").Value = 9.74
$ws.Range("Z5").Value = 6.34
$ws.Range("AA5").Value = 49.35
$ws.Range("AB5").Value = 10514.55
$ws.Range("AC5").Value = 4836
$ws.Range("AD5").Value = 11.29
$ws.Range("AE5").Value = 51341
$ws.Range("AF5").Value = 1.06
$ws.Range("AG5").Value = 400
$ws.Range("AH5").Value = 0.73
$ws.Range("AI5").Value = 8.27
$ws.Range("AJ5").Value = 123875069

$ws.Range("D6").Value = 67951
$ws.Range("E6").Value = 7027
$ws.Range("F6").Value = 7027
$ws.Range("G6").Value = 6977
$ws.Range("H6").Value = 5304
$ws.Range("I6").Value = 5222
$ws.Range("K6").Value = 97964
$ws.Range("L6").Value = 30061
$ws.Range("M6").Value = 67903
$ws.Range("N6").Value = 67655
$ws.Range("P6").Value = 619
$ws.Range("Q6").Value = 11235
$ws.Range("R6").Value = -5948
$ws.Range("S6").Value = -5828
$ws.Range("T6").Value = 3236
$ws.Range("U6").Value = 7999
$ws.Range("V6").Value = 16324
$ws.Range("W6").Value = 10.34
$ws.Range("X6").Value = 7.8
$ws.Range("Y6").Value = 7.96
$ws.Range("Z6").Value = 5.49
$ws.Range("AA6").Value = 44.27
$ws.Range("AB6").Value = 11251.23
$ws.Range("AC6").Value = 4216
$ws.Range("AD6").Value = 9.52
$ws.Range("AE6").Value = 54626
$ws.Range("AF6").Value = 0.74
$ws.Range("AG6").Value = 450
$ws.Range("AH6").Value = 1.12
$ws.Range("AI6").Value = 10.67
$ws.Range("AJ6").Value = 123875069

$ws.Range("D7").Value = 70046
$ws.Range("E7").Value = 5761
$ws.Range("G7").Value = 6515
$ws.Range("H7").Value = 5010
$ws.Range("I7").Value = 4897
$ws.Range("K7").Value = 104861
$ws.Range("L7").Value = 32138
$ws.Range("M7").Value = 72723
$ws.Range("N7").Value = 72408
$ws.Range("P7").Value = 620
$ws.Range("Q7").Value = 8869
$ws.Range("R7").Value = -3680
$ws.Range("S7").Value = -270
$ws.Range("T7").Value = 3821
$ws.Range("U7").Value = 5101
$ws.Range("W7").Value = 8.220000000000001
$ws.Range("X7").Value = 7.15
$ws.Range("Y7").Value = 6.99
$ws.Range("Z7").Value = 4.94
$ws.Range("AA7").Value = 44.19
$ws.Range("AC7").Value = 3953
$ws.Range("AD7").Value = 7.26
$ws.Range("AE7").Value = 58463
$ws.Range("AF7").Value = 0.49
$ws.Range("AG7").Value = 450
$ws.Range("AH7").Value = 1.57
$ws.Range("AI7").Value = 11.38

$ws.Range("D8").Value = 72069
$ws.Range("E8").Value = 6484
$ws.Range("G8").Value = 6938
$ws.Range("H8").Value = 5357
$ws.Range("I8").Value = 5261
$ws.Range("K8").Value = 108552
$ws.Range("L8").Value = 31162
$ws.Range("M8").Value = 77390
$ws.Range("N8").Value = 76989
$ws.Range("P8").Value = 620
$ws.Range("Q8").Value = 9786
$ws.Range("R8").Value = -5689
$ws.Range("S8").Value = -1438
$ws.Range("T8").Value = 4859
$ws.Range("U8").Value = 5025
$ws.Range("W8").Value = 9
$ws.Range("X8").Value = 7.43
$ws.Range("Y8").Value = 7.04
$ws.Range("Z8").Value = 5.02
$ws.Range("AA8").Value = 40.27
$ws.Range("AC8").Value = 4247
$ws.Range("AD8").Value = 6.76
$ws.Range("AE8").Value = 62162
$ws.Range("AF8").Value = 0.46
$ws.Range("AG8").Value = 463
$ws.Range("AH8").Value = 1.61
$ws.Range("AI8").Value = 10.9

$ws.Range("D9").Value = 74345
$ws.Range("E9").Value = 7216
$ws.Range("G9").Value = 7848
$ws.Range("H9").Value = 6040
$ws.Range("I9").Value = 5932
$ws.Range("K9").Value = 113606
$ws.Range("L9").Value = 30916
$ws.Range("M9").Value = 82690
$ws.Range("N9").Value = 82195
$ws.Range("P9").Value = 620
$ws.Range("Q9").Value = 10271
$ws.Range("R9").Value = -5944
$ws.Range("S9").Value = -1068
$ws.Range("T9").Value = 4935
$ws.Range("U9").Value = 6248
$ws.Range("W9").Value = 9.710000000000001
$ws.Range("X9").Value = 8.119999999999999
$ws.Range("Y9").Value = 7.45
$ws.Range("Z9").Value = 5.44
$ws.Range("AA9").Value = 37.39
$ws.Range("AC9").Value = 4789
$ws.Range("AD9").Value = 5.99
$ws.Range("AE9").Value = 66365
$ws.Range("AF9").Value = 0.43
$ws.Range("AG9").Value = 489
$ws.Range("AH9").Value = 1.7
$ws.Range("AI9").Value = 10.22
